$d = $word.ActiveDocument

$replacements = @(
    @("2025-06-02 Monday", "2025-06-03 Tuesday"),
    @("734×7=", "789×4="),
    @("743×6=", "362×3="),
    @("343×5=", "791×4="),
    @("623×5=", "201×9="),
    @("558×7=", "329×8="),
    @("401×4=", "262×7="),
    @("528×7=", "225×7="),
    @("457×8=", "443×3="),
    @("329×7=", "867×6="),
    @("234×6=", "553×3="),
    @("301×5=", "281×3="),
    @("826×8=", "425×8="),
    @("386×9=", "677×9="),
    @("702×7=", "288×2="),
    @("448×4=", "413×3="),
    @("713×6=", "847×7="),
    @("868×6=", "166×6="),
    @("710×6=", "199×8="),
    @("458×4=", "974×7="),
    @("911×9=", "993×6="),
    @("432×9=", "977×9="),
    @("227×8=", "406×9="),
    @("746×2=", "577×9="),
    @("257×3=", "239×4="),
    @("980×2=", "580×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
